$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1: was the shared string "销量" (a header label) -> now a date value (2020-08-15),
# formatted as a short date (built-in numFmtId 14).
$ws.Range("B1").Value = 44058
$ws.Range("B1").NumberFormat = "mm-dd-yy"

# Updated sales figures
$ws.Range("B2").Value = 2
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = 0

# Column B widened to fit the new date column
$ws.Columns("B").ColumnWidth = 9.8

# Update the active selection left over from editing
$ws.Range("J14").Select()
